$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 0. Drop every existing hyperlink. The backlog is being restructured
#    (rows 12-26 collapse into a smaller block plus a new column G
#    "overflow" list), so it's simplest to rebuild hyperlinks from
#    scratch in their final positions rather than try to patch them.
# -----------------------------------------------------------------
$ws.Hyperlinks.Delete()

# -----------------------------------------------------------------
# 1. Clear out the old stub rows (just an ID in col A + a link in col D)
#    for SM-8 and SM-10..SM-22. Using Clear() (not ClearContents) so
#    rows that end up fully empty drop out of sheetData/dimension.
# -----------------------------------------------------------------
$ws.Range("A12").Clear()
$ws.Range("D12").Clear()
$ws.Range("A14:A26").Clear()
$ws.Range("D14:D26").Clear()

# -----------------------------------------------------------------
# 2. New content for the Machine Learning / Spark cluster backlog
#    items that replace the removed stub rows.
# -----------------------------------------------------------------
function Set-Plain($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
}

Set-Plain "E12" "Nadja Ulrich, Matthias Fuhr (Data Scientists)"
Set-Plain "B12" "Machine Learning "

Set-Plain "B13" "Local Development with Spark (no Cloud)"

Set-Plain "B14" "Matching Amplab Reference Architecture"

Set-Plain "B15" "Real Time Word Cloud (Streaming)"

# -----------------------------------------------------------------
# 3. Rebuild all hyperlinks: the original D-column ones that remain,
#    plus the new G-column overflow list of issue links (8, 10-22) -
#    note issue 9 keeps its normal D13 slot since SM-9 survives.
# -----------------------------------------------------------------
function Add-Link($addr, $url, $text = $null) {
    if ($null -eq $text) { $text = $url }
    $r = $ws.Range($addr)
    $r.Value = $text
    $ws.Hyperlinks.Add($r, $url)
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
}

Add-Link "B2"  "https://github.com/Zuehlke/SHMACK/issues"
Add-Link "D5"  "https://github.com/Zuehlke/SHMACK/issues/1" "https://github.com/Zuehlke/SHMACK/issues/1 "
Add-Link "D6"  "https://github.com/Zuehlke/SHMACK/issues/2"
Add-Link "D11" "https://github.com/Zuehlke/SHMACK/issues/3"
Add-Link "D7"  "https://github.com/Zuehlke/SHMACK/issues/4"
Add-Link "D8"  "https://github.com/Zuehlke/SHMACK/issues/5"
Add-Link "D9"  "https://github.com/Zuehlke/SHMACK/issues/6"
Add-Link "D10" "https://github.com/Zuehlke/SHMACK/issues/7"
Add-Link "G11" "https://github.com/Zuehlke/SHMACK/issues/8"
Add-Link "D13" "https://github.com/Zuehlke/SHMACK/issues/9"
Add-Link "G13" "https://github.com/Zuehlke/SHMACK/issues/10"
Add-Link "G14" "https://github.com/Zuehlke/SHMACK/issues/11"
Add-Link "G15" "https://github.com/Zuehlke/SHMACK/issues/12"
Add-Link "G16" "https://github.com/Zuehlke/SHMACK/issues/13"
Add-Link "G17" "https://github.com/Zuehlke/SHMACK/issues/14"
Add-Link "G18" "https://github.com/Zuehlke/SHMACK/issues/15"
Add-Link "G19" "https://github.com/Zuehlke/SHMACK/issues/16"
Add-Link "G20" "https://github.com/Zuehlke/SHMACK/issues/17"
Add-Link "G21" "https://github.com/Zuehlke/SHMACK/issues/18"
Add-Link "G22" "https://github.com/Zuehlke/SHMACK/issues/19"
Add-Link "G23" "https://github.com/Zuehlke/SHMACK/issues/20"
Add-Link "G24" "https://github.com/Zuehlke/SHMACK/issues/21"
Add-Link "G25" "https://github.com/Zuehlke/SHMACK/issues/22"

# -----------------------------------------------------------------
# 4. Column G formatting (new "overflow" link column).
# -----------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 45.140625

# -----------------------------------------------------------------
# 5. Turn on the header AutoFilter (adds <autoFilter> + the hidden
#    _xlnm._FilterDatabase defined name).
# -----------------------------------------------------------------
$ws.Range("A4:E4").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Tabelle1!`$A`$4:`$E`$4")
$fdb.Visible = $false

# -----------------------------------------------------------------
# 6. Restore the selection to where the editor left off.
# -----------------------------------------------------------------
$ws.Range("E12").Select()
